$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.824.41"
$ws.Range("E2").Value = "  -1.66%  "
$ws.Range("D3").Value = "3.293.84"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("E4").Value = "  +0.01%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.60"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  -0.35%  "
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.67"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  -4.06%  "
$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.633"
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = "  +5.22%  "
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.125"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  -2.60%  "
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.66"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("E11").Value = "  -2.35%  "
$ws.Range("D12").Value = "3.867.22"
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("E13").Value = "  -3.48%  "
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.59"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  -2.82%  "
$ws.Range("D15").Value = "65.956.44"
$ws.Range("E15").Value = "  -1.88%  "
$ws.Range("D16").Value = "3.316.63"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("E17").Value = "  -1.72%  "
$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "435.87"
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = "  -1.50%  "
$ws.Range("E19").Value = "  -1.84%  "
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.29"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  -1.53%  "
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.42"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  -4.13%  "
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.41"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  -2.28%  "
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("B24").Value = "WrappedeETH"
$ws.Range("C24").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D24").Value = "3.445.50"
$ws.Range("E24").Value = "  +0.21%  "
$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.513"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("E26").Value = "  -4.08%  "
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.194"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  +3.70%  "
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.94"
$ws.Range("D28").Style = $origStyle
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("E30").Value = "  -1.42%  "
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.39"
$ws.Range("D31").Style = $origStyle
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("E33").Value = "  -3.37%  "
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.63"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  -2.11%  "
$ws.Range("E35").Value = "  -3.00%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "157.33"
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = "  -3.49%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.45"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  -4.72%  "
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "27.03"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  -1.39%  "
$ws.Range("E39").Value = "  -2.82%  "
$ws.Range("D40").Value = "2.780.36"
$ws.Range("E40").Value = "  +1.81%  "
$ws.Range("E41").Value = "  -0.39%  "
$ws.Range("E42").Value = "  -2.47%  "
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.29"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  +0.09%  "
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.08"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  -3.56%  "
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0658"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  -1.94%  "
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.29"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  -3.94%  "
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "320.48"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  -2.26%  "
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.49"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  -4.91%  "
$ws.Range("E49").Value = "  -1.18%  "
$ws.Range("E50").Value = "  +2.75%  "
$ws.Range("E51").Value = "  -0.04%  "
